# Generate Report for Handoff
# Regenerates the localization-status report with a new handoff file id
# (7a25559b-172d-473d-abed-ba2bce1ca9ce) and refreshed handoff timestamps,
# mirroring what the CI job produces each time it is re-run.

$wb = $excel.ActiveWorkbook

$oldGuid = "1774a987-5e0f-4ade-926c-03bba566df63"
$newGuid = "7a25559b-172d-473d-abed-ba2bce1ca9ce"

$oldHash = "94d34f1a519372d5350582d657069ba2d3642ee1"
$newHash = "535f2aaa5ab51a347d02b4811ea6d3228d3001d3"

$hyperlinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/56b3b5a63eace1e2598b1d64b8e53896b7b12790/e2e/$oldGuid.md"

$newMdName = "$newGuid.md"
$newMdPath = "e2e\$newGuid.md"

$newOverviewDate = "2016-09-02 03:06:56"
$newZhHandoffDate = "2016-09-02 03:06:52"
$newDeHandoffDate = "2016-09-02 03:06:56"

$newZhXlf = "$newGuid.$newHash.zh-cn.xlf"
$newDeXlf = "$newGuid.$newHash.de-de.xlf"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMdName
$wsOverview.Range("B2").Value = $newMdPath
$wsOverview.Range("G2").Value = $newOverviewDate

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $hyperlinkAddress, [Type]::Missing, [Type]::Missing, $newMdPath)

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = $newMdName
$wsZh.Range("G2").Value = $newZhXlf
$wsZh.Range("H2").Value = $newZhHandoffDate

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $hyperlinkAddress, [Type]::Missing, [Type]::Missing, $newMdName)

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = $newMdName
$wsDe.Range("G2").Value = $newDeXlf
$wsDe.Range("H2").Value = $newDeHandoffDate

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $hyperlinkAddress, [Type]::Missing, [Type]::Missing, $newMdName)
